# Update forecast error values for rows 2-15 (quarters Q9-Q22)
# to reflect corrected computation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2290253005182327
$ws.Range("C2").Value = 0.312612257312203
$ws.Range("D2").Value = 0.7824026871258996
$ws.Range("E2").Value = 0.4483004359774425
$ws.Range("F2").Value = 0.2817483719079676
$ws.Range("G2").Value = 0.2377561650759294
$ws.Range("H2").Value = 0.3001556178769585
$ws.Range("B3").Value = 0.08358695679397027
$ws.Range("C3").Value = 0.5533773866076669
$ws.Range("D3").Value = 0.2192751354592098
$ws.Range("E3").Value = 0.05272307138973492
$ws.Range("F3").Value = 0.00873086455769673
$ws.Range("G3").Value = 0.07113031735872585
$ws.Range("B4").Value = 0.4697904298136966
$ws.Range("C4").Value = 0.1356881786652395
$ws.Range("D4").Value = -0.03086388540423535
$ws.Range("E4").Value = -0.07485609223627354
$ws.Range("F4").Value = -0.01245663943524442
$ws.Range("G4").Value = -0.0456851157626226
$ws.Range("H4").Value = -0.03954651028770373
$ws.Range("I4").Value = -0.09447145597838819
$ws.Range("J4").Value = -0.09344167408882031
$ws.Range("B5").Value = -0.3341022511484572
$ws.Range("C5").Value = -0.500654315217932
$ws.Range("D5").Value = -0.5446465220499701
$ws.Range("E5").Value = -0.4822470692489411
$ws.Range("F5").Value = -0.5154755455763192
$ws.Range("G5").Value = -0.5093369401014003
$ws.Range("H5").Value = -0.5642618857920848
$ws.Range("I5").Value = -0.5632321039025169
$ws.Range("B6").Value = -0.1665520640694748
$ws.Range("C6").Value = -0.210544270901513
$ws.Range("D6").Value = -0.1481448181004839
$ws.Range("E6").Value = -0.1813732944278621
$ws.Range("F6").Value = -0.1752346889529432
$ws.Range("G6").Value = -0.2301596346436277
$ws.Range("H6").Value = -0.2291298527540598
$ws.Range("B7").Value = -0.04399220683203819
$ws.Range("C7").Value = 0.01840724596899092
$ws.Range("D7").Value = -0.01482123035838725
$ws.Range("E7").Value = -0.00868262488346838
$ws.Range("F7").Value = -0.06360757057415284
$ws.Range("G7").Value = -0.06257778868458495
$ws.Range("B8").Value = 0.06239945280102911
$ws.Range("C8").Value = 0.02917097647365094
$ws.Range("D8").Value = 0.03530958194856981
$ws.Range("E8").Value = -0.01961536374211465
$ws.Range("F8").Value = -0.01858558185254677
$ws.Range("G8").Value = -0.4422673433347777
$ws.Range("H8").Value = 0.0552385534021488
$ws.Range("I8").Value = -0.01630711523224489
$ws.Range("B9").Value = -0.03322847632737817
$ws.Range("C9").Value = -0.02708987085245931
$ws.Range("D9").Value = -0.08201481654314377
$ws.Range("E9").Value = -0.08098503465357589
$ws.Range("F9").Value = -0.5046667961358068
$ws.Range("G9").Value = -0.007160899398880316
$ws.Range("H9").Value = -0.078706568033274
$ws.Range("B10").Value = 0.006138605474918866
$ws.Range("C10").Value = -0.0487863402157656
$ws.Range("D10").Value = -0.04775655832619771
$ws.Range("E10").Value = -0.4714383198084287
$ws.Range("F10").Value = 0.02606757692849786
$ws.Range("G10").Value = -0.04547809170589583
$ws.Range("B11").Value = -0.05492494569068446
$ws.Range("C11").Value = -0.05389516380111658
$ws.Range("D11").Value = -0.4775769252833475
$ws.Range("E11").Value = 0.01992897145357899
$ws.Range("F11").Value = -0.0516166971808147
$ws.Range("B12").Value = 0.001029781889567885
$ws.Range("C12").Value = -0.422651979592663
$ws.Range("D12").Value = 0.07485391714426345
$ws.Range("E12").Value = 0.003308248509869764
$ws.Range("B13").Value = -0.4236817614822309
$ws.Range("C13").Value = 0.07382413525469557
$ws.Range("D13").Value = 0.002278466620301879
$ws.Range("B14").Value = 0.4975058967369265
$ws.Range("C14").Value = 0.4259602281025328
$ws.Range("B15").Value = -0.07154566863439368
